$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = "choice b"
$ws.Range("N2").Value = "Genus species"
$ws.Range("I3").Value = "choice b"
$ws.Range("N3").Value = "Genus species"

$ws.Range("N3").Select()
